# The underlying data (rows 2-22, weekly price records for Sandia /
# Mapocho Venta Directa de Santiago) was refreshed: the same 21 records
# were re-sorted / re-positioned in the sheet (a "weekly" resync per the
# commit message). Row 2 and row 13 keep their original content; every
# other row receives the content that used to live at a different row.
#
# Map of: new row number -> row number that currently holds the data
# that must end up there.
$rowMap = @{
    2  = 2
    3  = 12
    4  = 8
    5  = 9
    6  = 10
    7  = 11
    8  = 21
    9  = 22
    10 = 16
    11 = 17
    12 = 4
    13 = 13
    14 = 5
    15 = 6
    16 = 7
    17 = 3
    18 = 19
    19 = 20
    20 = 18
    21 = 15
    22 = 14
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 18  # R

# 1) Snapshot every source row's current values into a buffer first, so
#    that writes below never clobber a row we still need to read from
#    (the mapping contains cycles, e.g. 3 <- 12 <- 4 <- 8 <- 21 <- 15 <- 6 <- 10 <- 16 <- 7 <- 11 <- 17 <- 3).
$buffer = @{}
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    if (-not $buffer.ContainsKey($srcRow)) {
        $rowData = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowData[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $buffer[$srcRow] = $rowData
    }
}

# 2) Write the buffered data into its new home row.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $rowData = $buffer[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($newRow, $col).Value2 = $rowData[$col]
    }
}
